$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Week of 15/01 (row block 48-54) ---
$ws.Range("A48").Value = 46037
$ws.Range("B48").Value = "Implémentation des touches qui s'affichent lorsque l'utilisateur les touches."
$ws.Range("D48").Value = 1

$ws.Range("B49").Value = "Refactorisation du code."
$ws.Range("D49").Value = 1

$ws.Range("B50").Value = "Ajouter les boutons pour revenir vers la page pour choisir entre clavier et manette."
$ws.Range("D50").Value = 1

# --- Week of 16/01 (row block 55-61) ---
$ws.Range("A55").Value = 46038
$ws.Range("B55").Value = "Finalisation des tests"
$ws.Range("D55").Value = 0.5

$ws.Range("B56").Value = "Fait le diagramme de classe"
$ws.Range("D56").Value = 2.5

$ws.Range("B57").Value = "Fait le launch.json"
$ws.Range("D57").Value = 1.5

$ws.Range("B58").Value = "Fait la descente de code"
$ws.Range("D58").Value = 2.5

# --- Reflection texts (added after the task rows, matching shared-string order) ---
$ws.Range("B54").Value = "L’implémentation de l’affichage des touches a permis d’améliorer significativement l’interaction utilisateur. La refactorisation a rendu le code plus clair, structuré et plus facile à maintenir. L’ajout des boutons de navigation améliore l’ergonomie générale de l’application et facilite la transition entre les modes clavier et manette."
$ws.Rows.Item(54).RowHeight = 49.5

$ws.Range("B61").Value = "La finalisation des tests a permis de valider le bon fonctionnement global du projet. La réalisation du diagramme de classe a clarifié l’architecture et les relations entre les composants. La configuration du launch.json et la descente de code ont amélioré la compréhension du projet et facilité le débogage et la maintenance future."
$ws.Rows.Item(61).RowHeight = 36.75

# --- Update view / selection state to reflect scrolling down to the new entries ---
$excel.ActiveWindow.ScrollRow = 47
$ws.Range("J55").Select()
